$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current last row (145), shifting the
# existing last row down to 147.
$ws.Range("A145:R146").Insert()

# Row 145: new weekly record (Primera)
$ws.Cells.Item(145, 1).Value = 1
$ws.Cells.Item(145, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(145, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(145, 4).Value = 44890
$ws.Cells.Item(145, 5).Value = 15
$ws.Cells.Item(145, 6).Value = 100112042
$ws.Cells.Item(145, 7).Value = "Locoto"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Primera"
$ws.Cells.Item(145, 10).Value = 330
$ws.Cells.Item(145, 11).Value = 16000
$ws.Cells.Item(145, 12).Value = 17000
$ws.Cells.Item(145, 13).Value = 16545
$ws.Cells.Item(145, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(145, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(145, 16).Value = 827
$ws.Cells.Item(145, 17).Value = 20
$ws.Cells.Item(145, 18).Value = "Hortaliza"

# Row 146: new weekly record (Segunda)
$ws.Cells.Item(146, 1).Value = 1
$ws.Cells.Item(146, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(146, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(146, 4).Value = 44890
$ws.Cells.Item(146, 5).Value = 15
$ws.Cells.Item(146, 6).Value = 100112042
$ws.Cells.Item(146, 7).Value = "Locoto"
$ws.Cells.Item(146, 8).Value = "Sin especificar"
$ws.Cells.Item(146, 9).Value = "Segunda"
$ws.Cells.Item(146, 10).Value = 250
$ws.Cells.Item(146, 11).Value = 13000
$ws.Cells.Item(146, 12).Value = 14000
$ws.Cells.Item(146, 13).Value = 13400
$ws.Cells.Item(146, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(146, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(146, 16).Value = 670
$ws.Cells.Item(146, 17).Value = 20
$ws.Cells.Item(146, 18).Value = "Hortaliza"
